$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in row 4 (bug fix)
$ws.Range("E4").Value = 542
$ws.Range("F4").Value = 546
$ws.Range("G4").Value = 530
$ws.Range("H4").Value = 531
$ws.Range("I4").Value = 552
$ws.Range("J4").Value = 542
$ws.Range("K4").Value = 540

# Update the selection to A3 (as in the diff)
$ws.Range("A3").Select()
